$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words ("de", "del", "la", "las", "el",
#    "los", "y") wherever they appear as a standalone lowercase word inside
#    the state (column A) and municipality (column B) names.
#    Longer tokens are replaced before their shorter prefixes/substrings
#    (e.g. " del " before " de ", " las "/" los " before " la "/" lo ") so
#    that the wrong pattern never matches first.
$connectors = @(
  @(" del ", " Del "),
  @(" de ", " De "),
  @(" las ", " Las "),
  @(" la ", " La "),
  @(" los ", " Los "),
  @(" el ", " El "),
  @(" y ", " Y ")
)

foreach ($pair in $connectors) {
    $ws.Columns.Item(1).Replace($pair[0], $pair[1], 2) | Out-Null
    $ws.Columns.Item(2).Replace($pair[0], $pair[1], 2) | Out-Null
}

# 3. Row 179's state name had a stray "_x000D_" + newline artifact baked
#    into the text; rewrite it cleanly now that the wording is fixed.
$ws.Range("A179").Value = "Estado De México"

# 4. Tiny floating point recalculation drift on the Hidalgo state subtotal
#    row's percentage value.
$ws.Range("D451").Value = 0.09282106483398196

# 5. Drop the trailing footnote/metadata rows that are no longer part of
#    the tabular data, and let Excel shrink the sheet dimension to match.
$ws.Range("A1412:D1416").EntireRow.Delete()
